$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a handful of previously-computed values in row 9 ---
$ws.Range("A9").Value = -71.50511592146341
$ws.Range("B9").Value = -15.5254870799914
$ws.Range("C9").Value = 231276.4513187492
$ws.Range("D9").Value = 8281978.278665429
$ws.Range("G9").Value = 3829.822998046875

# --- Add the new "distancia" / "distancia_acumulada" columns ---
$ws.Range("H1").Value = "distancia"
$ws.Range("I1").Value = "distancia_acumulada"

# Match the look (bold, centered, bordered header) of the existing header row
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

$ws.Range("H3").Value = 72.34634270106764
$ws.Range("I3").Value = 72.34634270106764

$ws.Range("H4").Value = 45.80919622254626
$ws.Range("I4").Value = 118.1555389236139

$ws.Range("H5").Value = 47.94858565865728
$ws.Range("I5").Value = 166.1041245822712

$ws.Range("H6").Value = 53.98061808053652
$ws.Range("I6").Value = 220.0847426628077

$ws.Range("H7").Value = 59.08065155388351
$ws.Range("I7").Value = 279.1653942166912

$ws.Range("H8").Value = 37.15268647337687
$ws.Range("I8").Value = 316.3180806900681

$ws.Range("H9").Value = 33.21617207162281
$ws.Range("I9").Value = 349.5342527616909
